# Add a new 2D mandate variable "HMAN" (FTT-Heat minimum sales mandate)
# to the FTT-H variable list, and register it in the Time_Horizons sheet.

$wb = $excel.ActiveWorkbook

# --- 1. FTT-H sheet: insert the new HMAN variable row ------------------
$wsH = $wb.Worksheets.Item("FTT-H")

# New row sits right before the old "HETR" row (currently row 16).
$wsH.Rows.Item(16).Insert()

$wsH.Range("A16").Value = "HMAN"
$wsH.Range("B16").Value = 1
$wsH.Range("C16").Value = 33190000
$wsH.Range("D16").Value = "FTT-Heat minimum sales mandate"
$wsH.Range("E16").Value = "RSHORTTI"
$wsH.Range("F16").Value = "TIME"
$wsH.Range("G16").Value = 0
$wsH.Range("H16").Value = 0
$wsH.Range("I16").Value = "All"

# Highlight the dims cell with Excel's built-in "Neutral" cell style,
# flagging this as a 2-D variable like the others using that style.
$wsH.Range("C16").Style = "Neutral"

# --- 2. Time_Horizons sheet: register HMAN's time horizon --------------
$wsT = $wb.Worksheets.Item("Time_Horizons")

# New row sits right before the old row 31 (v110 / HWSA).
$wsT.Rows.Item(31).Insert()

$wsT.Range("A31").Value = "HMAN"
$wsT.Range("B31").Value = "tl_2001"

# Leave the Time_Horizons tab active/selected, matching the author's
# last-saved view state.
$wsT.Activate()
$wsT.Range("A32").Select()
